$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Certificados, cursos, badges" ---
$ws1 = $wb.Worksheets.Item("Certificados, cursos, badges")

# Row 76: IBM / governança de dados entry (previously blank except platform "IBM")
$ws1.Range("C76").Value = "`nFundamentos da governança de dados com o IBM Knowledge Catalog no IBM Cloud Pak for Data"
$ws1.Range("D76").Value = 3
$ws1.Range("E76").Value = 45577
$ws1.Range("F76").Value = "https://www.credly.com/badges/603ace26-de58-4b01-86a2-f57db8e29059/public_url"
$ws1.Range("G76").Value = "Ok"
$ws1.Range("H76").Value = "Ok"
$ws1.Range("I76").Value = 45577

# The leading line-break in C76 makes the host auto-expand the row; re-fit it
# back down to the sheet's normal (default) row height.
$ws1.Rows("76:76").AutoFit() | Out-Null

# Row 78: add the missing hyperlink in F78 (Acreditação Fundamentais Databricks)
$ws1.Range("F78").Value = "https://github.com/Phelipe-Sempreboni/certificates/blob/main/databricks-academy/fundamentals/acreditacao-fundamentais-databricks-portugues-br/certificado.pdf"

# Row 79: new Databricks Academy accreditation entry (was fully blank)
$ws1.Range("B79").Value = "Databricks"
$ws1.Range("C79").Value = "Acreditação da Academia - Fundamentos do Databricks"
$ws1.Range("D79").Value = 0
$ws1.Range("E79").Value = 45584
$ws1.Range("F79").Value = "https://credentials.databricks.com/4596b3d0-a718-480f-affb-385e2b40c6de"
$ws1.Range("G79").Value = "Ok"
$ws1.Range("H79").Value = "Ok"
$ws1.Range("I79").Value = 45584

# Register real hyperlinks for F76, F78, F79 (targets match the text above).
# Hyperlinks.Add changes the cell style/font size, so restore the original
# "hyperlink" cell look (font size 9, matching the rest of column F) afterwards.
$ws1.Hyperlinks.Add($ws1.Range("F78"), "https://github.com/Phelipe-Sempreboni/certificates/blob/main/databricks-academy/fundamentals/acreditacao-fundamentais-databricks-portugues-br/certificado.pdf") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F79"), "https://credentials.databricks.com/4596b3d0-a718-480f-affb-385e2b40c6de") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F76"), "https://www.credly.com/badges/603ace26-de58-4b01-86a2-f57db8e29059/public_url") | Out-Null

$ws1.Range("F76").Font.Size = 9
$ws1.Range("F78").Font.Size = 9
$ws1.Range("F79").Font.Size = 9

# Update the view: scrolled position / active selection
$ws1.Activate() | Out-Null
$ws1.Range("C80").Select() | Out-Null

# --- Sheet 2: "Formações" ---
$ws2 = $wb.Worksheets.Item("Formações")
$ws2.Activate() | Out-Null
$ws2.Range("B4").Select() | Out-Null

# Leave the first sheet as the active / tab-selected sheet, matching the source file.
$ws1.Activate() | Out-Null
